# Update vm_pu.xlsx results for "case with 380 kV" (Case_4_91).
# The bus voltage magnitude results in row 2..25 (columns B-F and I-N; column
# G stays pinned at the slack value of 1, column H has no data) are refreshed
# with the newly computed values from the 380 kV case.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that receive new values, in sheet order (G/H are left untouched).
$columns = @(2, 3, 4, 5, 6, 9, 10, 11, 12, 13, 14)   # B,C,D,E,F,I,J,K,L,M,N

# One row of new values per data row (sheet rows 2..25).
$newValues = @(
    @(1.02,1.054980614499226,1.061298055451379,1.068181144196707,1.074559391268119,1.047634003066163,1.059989869599541,1.064022496618958,1.070887014890782,1.077248283110919,1.061495176796909),
    @(1.02,1.056105691083205,1.062194530472678,1.069249192992797,1.075606001507586,1.047913896951509,1.060765593610772,1.064733358819364,1.071770365487352,1.078111478125053,1.062272002425213),
    @(1.02,1.056833840751527,1.06277443419167,1.069940695649012,1.076283446202158,1.048093338395268,1.061267113344771,1.065192522865719,1.072341753039275,1.078669632416698,1.06277423437476),
    @(1.02,1.057139991260612,1.063018183394105,1.070231499949235,1.076568295971339,1.048168375950903,1.061477850362543,1.065385361052565,1.07258191660274,1.07890418686091,1.062985270663269),
    @(1.02,1.057191397426798,1.063059107467522,1.070280332932942,1.076616126547272,1.048180951647729,1.061513228061037,1.065417728030007,1.072622238319829,1.078943564091065,1.063020698602153),
    @(1.02,1.056837931405995,1.06277769134711,1.06994458101087,1.076287252174498,1.048094342621224,1.061269929623429,1.065195100342974,1.072344962304296,1.078672766912669,1.062777054652857),
    @(1.02,1.055360808667362,1.06160105956094,1.068542012673788,1.074913052899707,1.047728939919444,1.060252117282131,1.06426290373402,1.071185588615852,1.077540084719883,1.061757796901289),
    @(1.02,1.052759054565187,1.059526351561621,1.066073601435097,1.072493212540329,1.047072280631854,1.058455341743202,1.062614046125304,1.069141098309517,1.075541166592108,1.059958469735002),
    @(1.02,1.051025272942602,1.058142330525983,1.064430072802862,1.070881115616892,1.046625927175618,1.05725529559908,1.061510638822347,1.067777077360804,1.074206548566483,1.058756719387708),
    @(1.02,1.05027468863476,1.057542826762052,1.063718898560835,1.070183327697751,1.046430617005897,1.056735139836062,1.061031862664471,1.067186196225627,1.07362816801777,1.058235824943845),
    @(1.02,1.049995910561144,1.057320112335027,1.063454809346773,1.069924177223712,1.046357764202477,1.056541851268934,1.060853874522439,1.06696667871935,1.073413259140686,1.05804226188478),
    @(1.02,1.050055708406408,1.057367886784702,1.063511454099626,1.069979764116229,1.046373405219175,1.056583315926882,1.060892060343491,1.067013767690272,1.073459361111638,1.058083785427298),
    @(1.02,1.050251644301132,1.057524417771332,1.063697067385042,1.070161905443366,1.046424601210907,1.0567191641672,1.061017153156576,1.067168051615423,1.073610405064638,1.0582198265877),
    @(1.02,1.050372369809369,1.057620857434487,1.063811439464889,1.070274133893897,1.04645610421507,1.056802854202029,1.061094207172826,1.067263105999075,1.073703458605798,1.058303635471977),
    @(1.02,1.051075089933999,1.058182113121478,1.06447728128372,1.070927430983712,1.046638846368671,1.057289805481089,1.06154239267887,1.067816286887014,1.07424492360756,1.058791278277708),
    @(1.02,1.051515928476404,1.058534116851721,1.064895075821304,1.071337296685991,1.046752930553519,1.057595115652703,1.061823261547899,1.068163215269014,1.0745844411091,1.059097022024785),
    @(1.02,1.051773077112169,1.058739414067687,1.065138815069245,1.071576389577651,1.04681927746671,1.057773146996153,1.061986991754119,1.068365548515887,1.074782429414842,1.059275306193163),
    @(1.02,1.05186076071754,1.058809411616723,1.065221931792615,1.071657918347007,1.046841866726137,1.057833842455208,1.062042803267761,1.068434534824868,1.074849930426557,1.059336087846732),
    @(1.02,1.051468629130261,1.058496352313569,1.064850245585992,1.071293319373296,1.04674071072103,1.057562364028699,1.061793136891475,1.068125995642434,1.074548018905933,1.059064223889718),
    @(1.02,1.050193945453464,1.057478324181231,1.063642406897471,1.070108268286026,1.046409533711977,1.056679162444539,1.060980320564812,1.067122619894135,1.073565928390204,1.058179768058004),
    @(1.02,1.049392630162859,1.056838064155958,1.062883410497019,1.06936340536772,1.046199539679904,1.056123398137715,1.060468406482503,1.066491538094391,1.072948029049129,1.057623214502118),
    @(1.02,1.049817410627716,1.05717749566363,1.063285729141037,1.069758249943579,1.046311029255474,1.056418063038767,1.060739863771807,1.06682610739201,1.0732756290492,1.057918297861128),
    @(1.02,1.051490001622228,1.05851341654633,1.064870502294952,1.071313190749994,1.046746232944283,1.057577163234948,1.061806749223031,1.068142813664203,1.074564476675299,1.059079044112538),
    @(1.02,1.053431540683224,1.060062869803824,1.066711378413663,1.073118600406161,1.047243555709897,1.058920238019683,1.063041050492993,1.069669829020623,1.076058288580056,1.06042402621752)
)

$startRow = 2
for ($i = 0; $i -lt $newValues.Count; $i++) {
    $rowValues = $newValues[$i]
    $rowNum = $startRow + $i
    for ($j = 0; $j -lt $columns.Count; $j++) {
        $ws.Cells.Item($rowNum, $columns[$j]).Value = $rowValues[$j]
    }
}
